$wb = $excel.ActiveWorkbook

# --- Rename "Requested quantity" header on "Weekly Quantity" sheet ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- Rename "Requested quantity" header on "Monthly Trend" sheet ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy header formatting (bold, centered, bordered) from an existing sheet header
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
        @(45144.99999999999, 61, -14.30917842252873, 136.3911490405453),
        @(45151.99999999999, 61, -11.94283942532171, 130.7144885150359),
        @(45186.99999999999, 62, -12.98149206693084, 133.9013913345954),
        @(45207.99999999999, 63, -6.998667681503075, 131.668654078593),
        @(45256.99999999999, 64, -5.122964196956632, 135.1696281652424),
        @(45263.99999999999, 64, -6.040328693296195, 136.6628711098233),
        @(45312.99999999999, 66, -8.547371341868498, 140.2141679405349),
        @(45382.99999999999, 67, -4.896583482531701, 135.6525061167258),
        @(45410.99999999999, 68, -0.249056127917128, 141.2368803149005),
        @(45445.99999999999, 69, -2.462834061675986, 146.0427704188151),
        @(45459.99999999999, 69, -0.2220982289409573, 140.9127209405634),
        @(45473.99999999999, 70, -3.273352168332071, 147.6290359624875),
        @(45480.99999999999, 70, -2.706256235176772, 139.9090859737202),
        @(45487.99999999999, 70, -4.95889321881076, 139.8365465623714),
        @(45501.99999999999, 70, -1.591471165791148, 142.7783604509277),
        @(45515.99999999999, 71, -6.386107834590546, 141.394420721497),
        @(45522.99999999999, 71, 2.476113889949145, 151.3140564381376),
        @(45536.99999999999, 71, 0.02481776719295787, 144.0295714717207),
        @(45543.99999999999, 72, 1.543965958890796, 149.348895280755),
        @(45564.99999999999, 72, 3.132625740798236, 138.7658753959599),
        @(45599.99999999999, 73, 4.429132147515068, 146.8190110075645),
        @(45620.99999999999, 74, 7.284550617081626, 150.099464106974),
        @(45627.99999999999, 74, 0.2056564534167741, 144.9671159183873),
        @(45634.99999999999, 74, 2.411889950728083, 146.4038228416947),
        @(45641.99999999999, 74, 3.283944022520237, 151.0578113378384),
        @(45648.99999999999, 74, 2.030602202318269, 144.110415105809),
        @(45655.99999999999, 74, -0.2121917675466934, 151.2206891667615),
        @(45662.99999999999, 75, 3.420810803196436, 142.4500848446792),
        @(45669.99999999999, 75, 0.3837964233432869, 144.950234480275),
        @(45676.99999999999, 75, 7.983006688082984, 146.7012267573515)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Copy the date-column style (numeric date format) from an existing sheet down column A
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
